$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 "Welcome to the Kiot" -> "Hulk" (new text, appended to sharedStrings)
$ws.Range("A2").Value = "Hulk"
# A3 "Batman" -> "Welcome to the Kiot" (reuses existing shared string)
$ws.Range("A3").Value = "Welcome to the Kiot"

# Insert a fresh row at position 4 carrying "Batman" (reuses existing shared string),
# pushing the old row4/row5 ("Welcome to the smartcliff" / "IronMan") down to 5/6.
$ws.Rows("4:4").Insert()
$ws.Range("A4").Value = "Batman"
# The insert copies the neighboring D-column formatting cell down; the target has no D4.
$ws.Range("D4").Clear()

# Add a new, empty, styled row 7 (style matches A3/A4/A5/A6 -> reuse that formatting).
$ws.Range("A3").Copy()
$ws.Range("A7").PasteSpecial(-4122) | Out-Null

# Restore / set explicit row heights to match the final layout.
$ws.Rows("2:2").RowHeight = 50.4
$ws.Rows("3:3").RowHeight = 39
$ws.Rows("4:4").RowHeight = 43.8
$ws.Rows("5:5").RowHeight = 37.8
$ws.Rows("6:6").RowHeight = 40.8
$ws.Rows("7:7").RowHeight = 54.6

# Widen column A (closest achievable width to the authored 42.21875).
$ws.Columns("A:A").ColumnWidth = 41.25

# Update the saved cursor/selection position.
$ws.Range("G8").Select() | Out-Null
